{"js": "// Updates the title date paragraph and overwrites every data cell in the\n// single results table, in row-major order, matching the target revision.\n// Run property formatting (fonts/sizes) is left untouched because only the\n// text contents of existing runs are modified.\n\nconst NEW_TITLE = \"2025-04-30 Wednesday\";\nconst NEW_GRID = [[\"93-35=58\", \"61-49=12\", \"31-2=29\", \"46-38=8\", \"96-37=59\"], [\"20-13=7\", \"41-35=6\", \"83-6=77\", \"58-19=39\", \"81-63=18\"], [\"26+66=92\", \"68+29=97\", \"18-9=9\", \"29+37=66\", \"72-46=26\"], [\"79+19=98\", \"87-19=68\", \"52-35=17\", \"9+16=25\", \"97-58=39\"], [\"44-36=8\", \"57+39=96\", \"75-58=17\", \"63-17=46\", \"59+24=83\"], [\"33-27=6\", \"90-34=56\", \"8+57=65\", \"19+64=83\", \"22+49=71\"], [\"7+67=74\", \"93-27=66\", \"53-19=34\", \"41-4=37\", \"5+7=12\"], [\"73-57=16\", \"29+65=94\", \"39+7=46\", \"84+7=91\", \"94-25=69\"], [\"7+67=74\", \"49+7=56\", \"95-38=57\", \"26+58=84\", \"93-25=68\"], [\"18-9=9\", \"3+29=32\", \"24-7=17\", \"33+38=71\", \"43-35=8\"], [\"38+39=77\", \"94-66=28\", \"70-63=7\", \"22-19=3\", \"61-13=48\"], [\"55+8=63\", \"90-83=7\", \"17+15=32\", \"29+5=34\", \"82-6=76\"], [\"5+67=72\", \"8+3=11\", \"86-37=49\", \"15+68=83\", \"61-42=19\"], [\"96-28=68\", \"28+29=57\", \"39+42=81\", \"74+18=92\", \"42-23=19\"], [\"41-34=7\", \"77-18=59\", \"3+29=32\", \"68+23=91\", \"27+19=46\"], [\"40-23=17\", \"13+29=42\", \"56-8=48\", \"39+59=98\", \"52+9=61\"], [\"38+43=81\", \"63-38=25\", \"98-39=59\", \"42-17=25\", \"70-25=45\"], [\"24+19=43\", \"29+24=53\", \"81-45=36\", \"70-56=14\", \"74-19=55\"], [\"94-15=79\", \"52-18=34\", \"24+7=31\", \"63-45=18\", \"8+58=66\"], [\"52+39=91\", \"70-44=26\", \"26+66=92\", \"15+9=24\", \"37+55=92\"]];\n\n// --- Title paragraph -------------------------------------------------\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length > 0) {\n  paragraphs.items[0].insertText(NEW_TITLE, Word.InsertLocation.replace);\n}\n\n// --- Table cells -------------------------------------------------------\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length > 0) {\n  const table = tables.items[0];\n  table.load(\"rowCount\");\n  await context.sync();\n\n  // Assigning `.values` rewrites each cell's text while Word keeps the\n  // existing paragraph/run formatting already present in that cell.\n  if (table.rowCount === NEW_GRID.length) {\n    table.values = NEW_GRID;\n  } else {\n    // Fallback: write cell by cell if the table shape is unexpected.\n    const rows = table.rows;\n    rows.load(\"items\");\n    await context.sync();\n    for (let r = 0; r < rows.items.length && r < NEW_GRID.length; r++) {\n      const cells = rows.items[r].cells;\n      cells.load(\"items\");\n      await context.sync();\n      for (let c = 0; c < cells.items.length && c < NEW_GRID[r].length; c++) {\n        cells.items[c].body.insertText(NEW_GRID[r][c], Word.InsertLocation.replace);\n      }\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Updates the title date paragraph and every data cell of the single\n# results table (20 rows x 5 columns), matching the target revision.\n# Only the text of each cell/paragraph is rewritten; Word preserves the\n# existing run formatting (fonts/sizes) already applied to that text.\n\n$NewTitle = \"2025-04-30 Wednesday\"\n\n$NewGrid = @(\n    @(\"93-35=58\", \"61-49=12\", \"31-2=29\", \"46-38=8\", \"96-37=59\"),\n    @(\"20-13=7\", \"41-35=6\", \"83-6=77\", \"58-19=39\", \"81-63=18\"),\n    @(\"26+66=92\", \"68+29=97\", \"18-9=9\", \"29+37=66\", \"72-46=26\"),\n    @(\"79+19=98\", \"87-19=68\", \"52-35=17\", \"9+16=25\", \"97-58=39\"),\n    @(\"44-36=8\", \"57+39=96\", \"75-58=17\", \"63-17=46\", \"59+24=83\"),\n    @(\"33-27=6\", \"90-34=56\", \"8+57=65\", \"19+64=83\", \"22+49=71\"),\n    @(\"7+67=74\", \"93-27=66\", \"53-19=34\", \"41-4=37\", \"5+7=12\"),\n    @(\"73-57=16\", \"29+65=94\", \"39+7=46\", \"84+7=91\", \"94-25=69\"),\n    @(\"7+67=74\", \"49+7=56\", \"95-38=57\", \"26+58=84\", \"93-25=68\"),\n    @(\"18-9=9\", \"3+29=32\", \"24-7=17\", \"33+38=71\", \"43-35=8\"),\n    @(\"38+39=77\", \"94-66=28\", \"70-63=7\", \"22-19=3\", \"61-13=48\"),\n    @(\"55+8=63\", \"90-83=7\", \"17+15=32\", \"29+5=34\", \"82-6=76\"),\n    @(\"5+67=72\", \"8+3=11\", \"86-37=49\", \"15+68=83\", \"61-42=19\"),\n    @(\"96-28=68\", \"28+29=57\", \"39+42=81\", \"74+18=92\", \"42-23=19\"),\n    @(\"41-34=7\", \"77-18=59\", \"3+29=32\", \"68+23=91\", \"27+19=46\"),\n    @(\"40-23=17\", \"13+29=42\", \"56-8=48\", \"39+59=98\", \"52+9=61\"),\n    @(\"38+43=81\", \"63-38=25\", \"98-39=59\", \"42-17=25\", \"70-25=45\"),\n    @(\"24+19=43\", \"29+24=53\", \"81-45=36\", \"70-56=14\", \"74-19=55\"),\n    @(\"94-15=79\", \"52-18=34\", \"24+7=31\", \"63-45=18\", \"8+58=66\"),\n    @(\"52+39=91\", \"70-44=26\", \"26+66=92\", \"15+9=24\", \"37+55=92\"),\n)\n\n$d = $word.ActiveDocument\n\n# --- Title paragraph ---------------------------------------------------\n$d.Paragraphs.Item(1).Range.Text = $NewTitle\n\n# --- Table cells ---------------------------------------------------------\n$table = $d.Tables.Item(1)\n$rowCount = [Math]::Min($table.Rows.Count, $NewGrid.Count)\n$colCount = [Math]::Min($table.Columns.Count, $NewGrid[0].Count)\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $table.Cell($r, $c).Range.Text = $NewGrid[$r - 1][$c - 1]\n    }\n}\n\n"}
